$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying the existing "2022-Q2"
#    sheet (so all styles / column layout / header row are inherited
#    identically), placed right after "总计".
# ------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy($null, $summarySheet)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template had 8 data rows (rows 2-8); the new sheet only needs 2.
$q3.Rows("4:8").Delete()

# Row 2: fund 011243
$q3.Cells.Item(2, 2).Value = "'011243"
$q3.Cells.Item(2, 2).Style = "Normal"
$q3.Cells.Item(2, 3).Value = "万家惠裕回报6个月持有期混合A"
$q3.Cells.Item(2, 4).Value = "'1.54"
$q3.Cells.Item(2, 4).Style = "Normal"
$q3.Cells.Item(2, 5).Value = "'27.67"
$q3.Cells.Item(2, 5).Style = "Normal"
$q3.Cells.Item(2, 6).Value = "'0.68"
$q3.Cells.Item(2, 6).Style = "Normal"
$q3.Cells.Item(2, 7).Value = "'0.0105"
$q3.Cells.Item(2, 7).Style = "Normal"
$q3.Cells.Item(2, 8).Value = 9

# Row 3: fund 011244
$q3.Cells.Item(3, 2).Value = "'011244"
$q3.Cells.Item(3, 2).Style = "Normal"
$q3.Cells.Item(3, 3).Value = "万家惠裕回报6个月持有期混合C"
$q3.Cells.Item(3, 4).Value = "'0.12"
$q3.Cells.Item(3, 4).Style = "Normal"
$q3.Cells.Item(3, 5).Value = "'27.67"
$q3.Cells.Item(3, 5).Style = "Normal"
$q3.Cells.Item(3, 6).Value = "'0.68"
$q3.Cells.Item(3, 6).Style = "Normal"
$q3.Cells.Item(3, 7).Value = "'0.0008"
$q3.Cells.Item(3, 7).Style = "Normal"
$q3.Cells.Item(3, 8).Value = 9

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: prepend a 2022-Q3 row, pushing
#    all the other quarters down by one row, and append the trailing
#    index row.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$data = @(
    @(0, "2022-Q3", 2, 0.01),
    @(1, "2022-Q2", 7, 0.33),
    @(2, "2022-Q1", 4, 0.2),
    @(3, "2021-Q4", 7, 0.49),
    @(4, "2021-Q3", 13, 1.64),
    @(5, "2021-Q2", 24, 2.16),
    @(6, "2021-Q1", 35, 8.76),
    @(7, "2020-Q4", 31, 24.86)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $summary.Cells.Item($row, 1).Value = $data[$i][0]
    $summary.Cells.Item($row, 2).Value = $data[$i][1]
    $summary.Cells.Item($row, 3).Value = $data[$i][2]
    $summary.Cells.Item($row, 4).Value = $data[$i][3]
}

# Row 9's "A" cell is brand new territory (old sheet only went to row 8)
# so it has no inherited style yet - copy the style from row 8's A cell.
$summary.Cells.Item(8, 1).Copy()
$summary.Cells.Item(9, 1).PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Restore "总计" as the active sheet/selection.
# ------------------------------------------------------------------
$summary.Activate()
$summary.Range("A1").Select() | Out-Null
